$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  5"
$ws.Range("C9").Value = "Report Covering the Week  1/30/2023  Through  2/5/2023"

# --- Fix cell styles that change number<->text representation (paste formats from stable donor cells) ---
$ws.Range("H14").Copy() | Out-Null
$ws.Range("M14").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null
$ws.Range("A15").Copy() | Out-Null
$ws.Range("C26").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("F22").PasteSpecial(-4122) | Out-Null
$ws.Range("D23").Copy() | Out-Null
$ws.Range("I22").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Cell value updates ---
$ws.Range("M14").Value = -100
$ws.Range("C15").Value = "0"
$ws.Range("E15").Value = -100
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -33.333333333333
$ws.Range("J15").Value = 3
$ws.Range("K15").Value = -33.333333333333
$ws.Range("M15").Value = -33.333333333333
$ws.Range("N15").Value = -77.777777777777
$ws.Range("D16").Value = 3
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = -47.368421052631
$ws.Range("J16").Value = 22
$ws.Range("K16").Value = -31.818181818181
$ws.Range("L16").Value = 66.666666666666
$ws.Range("M16").Value = -44.444444444444
$ws.Range("N16").Value = -85.981308411215
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 33
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 42
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = 5
$ws.Range("L17").Value = 16.666666666666
$ws.Range("M17").Value = 35.483870967741
$ws.Range("N17").Value = -42.465753424657
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -28.571428571428
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 17
$ws.Range("H18").Value = -29.411764705882
$ws.Range("I18").Value = 23
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = 4.545454545454
$ws.Range("L18").Value = 4.545454545454
$ws.Range("M18").Value = 53.333333333333
$ws.Range("N18").Value = -74.157303370786
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = 21
$ws.Range("H19").Value = -14.285714285714
$ws.Range("I19").Value = 22
$ws.Range("J19").Value = 31
$ws.Range("K19").Value = -29.032258064516
$ws.Range("L19").Value = -21.428571428571
$ws.Range("M19").Value = -15.384615384615
$ws.Range("N19").Value = -56
$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = -16.666666666666
$ws.Range("I20").Value = 6
$ws.Range("K20").Value = -14.285714285714
$ws.Range("L20").Value = -25
$ws.Range("M20").Value = 20
$ws.Range("N20").Value = -80
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -9.523809523809
$ws.Range("F21").Value = 80
$ws.Range("H21").Value = -21.56862745098
$ws.Range("I21").Value = 110
$ws.Range("J21").Value = 128
$ws.Range("K21").Value = -14.0625
$ws.Range("L21").Value = 5.76923076923
$ws.Range("M21").Value = 1.851851851851
$ws.Range("N21").Value = -69.696969696969
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 1
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -50
$ws.Range("M22").Value = -50
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 400
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = 46.153846153846
$ws.Range("I23").Value = 22
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = 29.411764705882
$ws.Range("L23").Value = 46.666666666666
$ws.Range("M23").Value = 69.230769230769
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 6.25
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 7.692307692307
$ws.Range("I24").Value = 86
$ws.Range("J24").Value = 82
$ws.Range("K24").Value = 4.878048780487
$ws.Range("L24").Value = 4.878048780487
$ws.Range("M24").Value = 32.307692307692
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 32
$ws.Range("H25").Value = 59.375
$ws.Range("I25").Value = 60
$ws.Range("J25").Value = 48
$ws.Range("K25").Value = 25
$ws.Range("L25").Value = 62.162162162162
$ws.Range("M25").Value = -30.232558139534
$ws.Range("C26").Value = "0"
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 3
$ws.Range("K26").Value = 0
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -80
$ws.Range("I27").Value = 3
$ws.Range("J27").Value = 6
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = -75
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -25
$ws.Range("N28").Value = -78.571428571428
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -81.818181818181
